$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.05354133333333333
$ws.Range("H2").Value = 0.160624
$ws.Range("I2").Value = 0.00209946492164722
$ws.Range("J2").Value = 0.00209946492164722
$ws.Range("M2").Value = 1.845768666666667
$ws.Range("N2").Value = 5.537306
$ws.Range("O2").Value = 0.01459089321241885
$ws.Range("P2").Value = 0.01459089321241885
$ws.Range("Q2").Value = 0.09882491543822221
$ws.Range("R2").Value = 0.8894242389439999
$ws.Range("S2").Value = 0.00003063306847497389
$ws.Range("T2").Value = 0.00003063306847497389
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.05354133333333333
$ws.Range("H3").Value = 0.160624
$ws.Range("I3").Value = 0.00209946492164722
$ws.Range("J3").Value = 0.00209946492164722
$ws.Range("O3").Value = 0.6557810310272387
$ws.Range("P3").Value = 0.6557810310272387
$ws.Range("Q3").Value = 4.441640685992889
$ws.Range("R3").Value = 39.974766173936
$ws.Range("S3").Value = 0.001376789270923335
$ws.Range("T3").Value = 0.001376789270923335
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.05354133333333333
$ws.Range("H4").Value = 0.160624
$ws.Range("I4").Value = 0.00209946492164722
$ws.Range("J4").Value = 0.00209946492164722
$ws.Range("M4").Value = 41.69841866666667
$ws.Range("N4").Value = 125.095256
$ws.Range("O4").Value = 0.3296280757603424
$ws.Range("P4").Value = 0.3296280757603424
$ws.Range("Q4").Value = 2.232588933304889
$ws.Range("R4").Value = 20.093300399744
$ws.Range("S4").Value = 0.0006920425822489111
$ws.Range("T4").Value = 0.0006920425822489111
$ws.Range("I5").Value = 0.05460670042535784
$ws.Range("J5").Value = 0.05460670042535784
$ws.Range("M5").Value = 1.845768666666667
$ws.Range("N5").Value = 5.537306
$ws.Range("O5").Value = 0.01459089321241885
$ws.Range("P5").Value = 0.01459089321241885
$ws.Range("Q5").Value = 2.570418060456222
$ws.Range("R5").Value = 23.133762544106
$ws.Range("S5").Value = 0.0007967605345889433
$ws.Range("T5").Value = 0.0007967605345889433
$ws.Range("I6").Value = 0.05460670042535784
$ws.Range("J6").Value = 0.05460670042535784
$ws.Range("O6").Value = 0.6557810310272387
$ws.Range("P6").Value = 0.6557810310272387
$ws.Range("S6").Value = 0.03581003830593672
$ws.Range("T6").Value = 0.03581003830593672
$ws.Range("I7").Value = 0.05460670042535784
$ws.Range("J7").Value = 0.05460670042535784
$ws.Range("M7").Value = 41.69841866666667
$ws.Range("N7").Value = 125.095256
$ws.Range("O7").Value = 0.3296280757603424
$ws.Range("P7").Value = 0.3296280757603424
$ws.Range("Q7").Value = 58.0692317346729
$ws.Range("R7").Value = 522.6230856120561
$ws.Range("S7").Value = 0.01799990158483218
$ws.Range("T7").Value = 0.01799990158483218
$ws.Range("G8").Value = 24.05622933333333
$ws.Range("H8").Value = 72.168688
$ws.Range("I8").Value = 0.943293834652995
$ws.Range("J8").Value = 0.943293834652995
$ws.Range("M8").Value = 1.845768666666667
$ws.Range("N8").Value = 5.537306
$ws.Range("O8").Value = 0.01459089321241885
$ws.Range("P8").Value = 0.01459089321241885
$ws.Range("Q8").Value = 44.40223434161422
$ws.Range("R8").Value = 399.620109074528
$ws.Range("S8").Value = 0.01376349960935493
$ws.Range("T8").Value = 0.01376349960935493
$ws.Range("G9").Value = 24.05622933333333
$ws.Range("H9").Value = 72.168688
$ws.Range("I9").Value = 0.943293834652995
$ws.Range("J9").Value = 0.943293834652995
$ws.Range("O9").Value = 0.6557810310272387
$ws.Range("P9").Value = 0.6557810310272387
$ws.Range("Q9").Value = 1995.638141719337
$ws.Range("R9").Value = 17960.74327547403
$ws.Range("S9").Value = 0.6185942034503787
$ws.Range("T9").Value = 0.6185942034503787
$ws.Range("G10").Value = 24.05622933333333
$ws.Range("H10").Value = 72.168688
$ws.Range("I10").Value = 0.943293834652995
$ws.Range("J10").Value = 0.943293834652995
$ws.Range("M10").Value = 41.69841866666667
$ws.Range("N10").Value = 125.095256
$ws.Range("O10").Value = 0.3296280757603424
$ws.Range("P10").Value = 0.3296280757603424
$ws.Range("Q10").Value = 1003.106722282681
$ws.Range("R10").Value = 9027.960500544128
$ws.Range("S10").Value = 0.3109361315932613
$ws.Range("T10").Value = 0.3109361315932613
